$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update date text (slash -> dash). A leading apostrophe forces these to stay
# plain text cells instead of being auto-parsed into date serial numbers by
# Excel's locale-aware value parser (important for day <= 12 values, which
# are otherwise ambiguous with MM-DD-YYYY).
$ws.Range("A3").Value  = "'28-07-2022"
$ws.Range("A4").Value  = "'01-08-2022"
$ws.Range("A5").Value  = "'04-08-2022"
$ws.Range("A6").Value  = "'08-08-2022"
$ws.Range("A7").Value  = "'11-08-2022"
$ws.Range("A8").Value  = "'15-08-2022"
$ws.Range("A9").Value  = "'18-08-2022"
$ws.Range("A10").Value = "'22-08-2022"
$ws.Range("A11").Value = "'25-08-2022"
$ws.Range("A12").Value = "'29-08-2022"
$ws.Range("A13").Value = "'01-09-2022"
$ws.Range("A14").Value = "'05-09-2022"
$ws.Range("A15").Value = "'08-09-2022"
$ws.Range("A16").Value = "'12-09-2022"
$ws.Range("A17").Value = "'15-09-2022"
$ws.Range("A18").Value = "'19-09-2022"
$ws.Range("A19").Value = "'22-09-2022"
$ws.Range("A20").Value = "'26-09-2022"
$ws.Range("A21").Value = "'29-09-2022"

# Updated attendance counters
$ws.Range("D3").Value  = 1
$ws.Range("G3").Value  = 1

$ws.Range("D4").Value  = 1
$ws.Range("E4").Value  = 1
$ws.Range("H4").Value  = 0

$ws.Range("D5").Value  = 1
$ws.Range("E5").Value  = 1
$ws.Range("H5").Value  = 0

$ws.Range("D11").Value = 1
$ws.Range("E11").Value = 1
$ws.Range("H11").Value = 0

$ws.Range("D12").Value = 1
$ws.Range("E12").Value = 1
$ws.Range("H12").Value = 0

$ws.Range("D16").Value = 1
$ws.Range("E16").Value = 1
$ws.Range("H16").Value = 0
